$d = $word.ActiveDocument

# Locate the run containing "do Cục CSQLHC về TTXH cấp"
$r = $d.Content
$found = $r.Find.Execute("do Cục CSQLHC về TTXH cấp", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $r.Start
    $end = $r.End

    # "do " = 3 chars, "Cục CSQLHC về TTXH" = 18 chars, " cấp" = 4 chars
    $midStart = $start + 3
    $midEnd = $end - 4

    $mid = $d.Range($midStart, $midEnd)

    # Toggling the (identical) character formatting on just the middle
    # sub-range forces Word to split the single run into three runs
    # ("do " / "Cục CSQLHC về TTXH" / " cấp") while preserving the
    # original formatting (black font color) on all of them.
    $mid.Font.Color = 1
    $mid.Font.Color = 0
}
